$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in this week's Roboflow annotation report row (7/11/2025)
$ws.Range("D56").Value = (Get-Date -Year 2025 -Month 11 -Day 7 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("E56").Value = 150
$ws.Range("F56").Value = 776
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 1012
$ws.Range("J56").Value = "N/A"

# Scroll the view down to the newly-entered row and move the selection to the next row
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Application.ActiveWindow.ScrollRow = 34
$ws.Range("D57").Select()
